$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mega 2560")

# Pre-size row 19 (45pt) before the wrapped text goes in it. The picture
# floating below (anchored row 11 -> row 47) keeps a fixed pixel size, so
# growing row 19 shifts its lower anchor row up; restore the picture's
# height right after so only the anchor (not the displayed size) moves.
$shp = $ws.Shapes.Item("Grafik 1")
$origHeight = $shp.Height
$ws.Rows.Item(19).RowHeight = 45
$shp.Height = $origHeight

# New PIN-assignment rows for the coin-acceptor power switch, cashless
# (onyx) payment and the NV10 bill reader, plus the two new interrupt pins
# feeding the NV10's channel-open lines.
# Cell values are written in the same order the strings first appear in
# the target shared-string table so the generated sharedStrings.xml lines
# up with the authored edit.
$ws.Range("E20").Value = "cashless payment onyx"
$ws.Range("C19").Value = "coin power"
$ws.Range("E19").Value = "coin acceptor power - allows to turn off coin acceptor when all compartments are`nempty"
$ws.Range("E19").WrapText = $true
$ws.Range("E21").Value = "NV10 bill reader"
$ws.Range("E11").Value = "nv 10 channel 1 open"
$ws.Range("E12").Value = "nv 10 channel 2 open"
$ws.Range("E13").Value = "onyx inhibit"
$ws.Range("C20").Value = "Interrupt"
$ws.Range("C21").Value = "Interrupt"

$ws.Activate() | Out-Null
$ws.Range("E16").Select() | Out-Null
